$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (row 1): I1 = "I0", J1 = "IF"
# Match the style used by the other header cells (bold, centered, bordered)
# by copying the formatting from H1 before setting the new values.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data rows 2-19 for columns I (I0) and J (IF)
$values = @{
    2  = @(6, 6)
    3  = @(7, 7)
    4  = @(13, 13)
    5  = @(8, 8)
    6  = @(4, 5)
    7  = @(9, 9)
    8  = @(5, 6)
    9  = @(9, 9)
    10 = @(8, 8)
    11 = @(2, 3)
    12 = @(8, 9)
    13 = @(8, 8)
    14 = @(9, 9)
    15 = @(4, 4)
    16 = @(8, 8)
    17 = @(8, 8)
    18 = @(9, 9)
    19 = @(7, 7)
}

foreach ($row in $values.Keys) {
    $pair = $values[$row]
    $ws.Cells.Item($row, 9).Value = $pair[0]
    $ws.Cells.Item($row, 10).Value = $pair[1]
}
